# Apply updated crypto price/volume data as captured by the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.557.27"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.882.96"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'246.40"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.4732"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.2893"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "'0.06538"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "'22.26"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").Value = "'99.70"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").Value = "'0.7606"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("D13").Value = "'0.07826"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "1.880.98"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "'5.235"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "'284.24"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "30.543.49"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'13.18"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'0.000007521"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'0.9985"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "2.127.34"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'5.344"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'6.424"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "'9.171"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'163.50"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "'19.01"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'1.906"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "'0.09766"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'1.328"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Value = "'4.250"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "'4.186"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'0.6991"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'2.763"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").Value = "'0.01905"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'2.875"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "'6.310"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'75.32"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'1.974"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "'0.4248"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8372"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "'9.948"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").Value = "'101.43"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'7.019"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'35.25"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'0.05794"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'0.3956"
$ws.Range("E51").Value = "  -0.17%  "
